$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @($true,  "hyper_heuristic", 3888.081174736501),
    @($false, "hyper_heuristic", 5045.427369818298),
    @($true,  "simple",          3652.688798657419),
    @($false, "simple",          3644.774037914917)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
